$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 changes from "DDP-001-A2" to "DDP-001-A1"
$ws.Range("A2").Value = "DDP-001-A1"

# Rows 3-11 gain columns I:Y, populated with the same values already
# present in I2:Y2 (the template/header row for that data).
$src = $ws.Range("I2:Y2")
for ($r = 3; $r -le 11; $r++) {
    $dst = $ws.Range("I" + $r + ":Y" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# Move/record the active selection at A2
$ws.Range("A2").Select()
